$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1000207.4
$ws.Range("I6").Value = 1250034.2
$ws.Range("K6").Value = 3750102.6
$ws.Range("M6").Value = -3749990.6
$ws.Range("H51").Value = 8646.467000000001
$ws.Range("I51").Value = 27050
$ws.Range("J51").Value = 1954.2727
$ws.Range("K51").Value = 27050
$ws.Range("L51").Value = 1954.2727
$ws.Range("M51").Value = -26566
$ws.Range("N51").Value = -2922.2727
$ws.Range("H64").Value = 85822.414
$ws.Range("I64").Value = 335033
$ws.Range("J64").Value = 2752.2222
$ws.Range("K64").Value = 335033
$ws.Range("L64").Value = 2752.2222
$ws.Range("M64").Value = -334785
$ws.Range("N64").Value = -3248.2222
$ws.Range("H67").Value = 85822.414
$ws.Range("I67").Value = 335033
$ws.Range("J67").Value = 2752.2222
$ws.Range("K67").Value = 335033
$ws.Range("L67").Value = 2752.2222
$ws.Range("M67").Value = -334175
$ws.Range("N67").Value = -4468.2222
$ws.Range("H100").Value = 774.1111
$ws.Range("I100").Value = 630.3333
$ws.Range("J100").Value = 917.8889
$ws.Range("K100").Value = 630.3333
$ws.Range("L100").Value = 917.8889
$ws.Range("M100").Value = -89.33330000000001
$ws.Range("N100").Value = -1999.8889
$ws.Range("H112").Value = 964.5862
$ws.Range("J112").Value = 987.8889
$ws.Range("L112").Value = 2963.6667
$ws.Range("N112").Value = -5179.6667
$ws.Range("H138").Value = 1187.5741
$ws.Range("I138").Value = 998.5306399999999
$ws.Range("J138").Value = 3040.2
$ws.Range("K138").Value = 2995.59192
$ws.Range("L138").Value = 9120.599999999999
$ws.Range("M138").Value = 2144.40808
$ws.Range("N138").Value = -19400.6
$ws.Range("H141").Value = 2835.1936
$ws.Range("I141").Value = 2836.7083
$ws.Range("J141").Value = 2830
$ws.Range("K141").Value = 8510.124899999999
$ws.Range("L141").Value = 8490
$ws.Range("M141").Value = -3330.124899999999
$ws.Range("N141").Value = -18850
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2144.634
$ws.Range("I61").Value = 1072.1177
$ws.Range("K61").Value = 1072.1177
$ws.Range("M61").Value = -860.1177
$ws.Range("H132").Value = 2465.611
$ws.Range("I132").Value = 2048.0833
$ws.Range("J132").Value = 3300.6667
$ws.Range("K132").Value = 6144.249899999999
$ws.Range("L132").Value = 9902.000100000001
$ws.Range("M132").Value = -3614.249899999999
$ws.Range("N132").Value = -14962.0001
$ws.Range("H136").Value = 2144.634
$ws.Range("I136").Value = 1072.1177
$ws.Range("K136").Value = 3216.3531
$ws.Range("M136").Value = -666.3531000000003
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 65449.75
$ws.Range("I20").Value = 85499.75
$ws.Range("J20").Value = 5299.75
$ws.Range("K20").Value = 85499.75
$ws.Range("L20").Value = 5299.75
$ws.Range("M20").Value = -85252.75
$ws.Range("N20").Value = -5793.75
$ws.Range("H134").Value = 4884.8184
$ws.Range("I134").Value = 5012.9546
$ws.Range("J134").Value = 4628.5454
$ws.Range("K134").Value = 15038.8638
$ws.Range("L134").Value = 13885.6362
$ws.Range("M134").Value = -12503.8638
$ws.Range("N134").Value = -18955.6362
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1075.5
$ws.Range("I16").Value = 1122.2
$ws.Range("J16").Value = 997.6667
$ws.Range("K16").Value = 1122.2
$ws.Range("L16").Value = 997.6667
$ws.Range("M16").Value = -835.2
$ws.Range("N16").Value = -1571.6667
$ws.Range("H31").Value = 26298.465
$ws.Range("I31").Value = 1260.9565
$ws.Range("J31").Value = 42751.688
$ws.Range("K31").Value = 1260.9565
$ws.Range("L31").Value = 42751.688
$ws.Range("M31").Value = -965.9565
$ws.Range("N31").Value = -43341.688
$ws.Range("H34").Value = 26298.465
$ws.Range("I34").Value = 1260.9565
$ws.Range("J34").Value = 42751.688
$ws.Range("K34").Value = 1260.9565
$ws.Range("L34").Value = 42751.688
$ws.Range("M34").Value = -1058.9565
$ws.Range("N34").Value = -43155.688
$ws.Range("H58").Value = 1132.2162
$ws.Range("I58").Value = 1108.3889
$ws.Range("J58").Value = 1990
$ws.Range("K58").Value = 1108.3889
$ws.Range("L58").Value = 1990
$ws.Range("M58").Value = -905.3888999999999
$ws.Range("N58").Value = -2396
$ws.Range("H107").Value = 780.8421
$ws.Range("I107").Value = 912.0833
$ws.Range("K107").Value = 912.0833
$ws.Range("M107").Value = 1007.9167
$ws.Range("H113").Value = 1075.5
$ws.Range("I113").Value = 1122.2
$ws.Range("J113").Value = 997.6667
$ws.Range("K113").Value = 1122.2
$ws.Range("L113").Value = 997.6667
$ws.Range("M113").Value = 1047.8
$ws.Range("N113").Value = -5337.6667
$ws.Range("H132").Value = 23439668
$ws.Range("I132").Value = 22224304
$ws.Range("J132").Value = 26318164
$ws.Range("K132").Value = 66672912
$ws.Range("L132").Value = 78954492
$ws.Range("M132").Value = -66670382
$ws.Range("N132").Value = -78959552
$ws.Range("H136").Value = 1132.2162
$ws.Range("I136").Value = 1108.3889
$ws.Range("J136").Value = 1990
$ws.Range("K136").Value = 3325.1667
$ws.Range("L136").Value = 5970
$ws.Range("M136").Value = -775.1666999999998
$ws.Range("N136").Value = -11070
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 114.94118
$ws.Range("I2").Value = 172.28572
$ws.Range("J2").Value = 74.8
$ws.Range("K2").Value = 1033.71432
$ws.Range("L2").Value = 448.8
$ws.Range("M2").Value = -920.71432
$ws.Range("N2").Value = -674.8
$ws.Range("H38").Value = 40.81818
$ws.Range("J38").Value = 64.833336
$ws.Range("L38").Value = 194.500008
$ws.Range("N38").Value = -888.500008
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 673741.0600000001
$ws.Range("I107").Value = 381.08334
$ws.Range("J107").Value = 3367181
$ws.Range("K107").Value = 381.08334
$ws.Range("L107").Value = 3367181
$ws.Range("M107").Value = 1538.91666
$ws.Range("N107").Value = -3371021
$ws.Range("H132").Value = 2917.0908
$ws.Range("I132").Value = 2194.4285
$ws.Range("J132").Value = 4181.75
$ws.Range("K132").Value = 6583.2855
$ws.Range("L132").Value = 12545.25
$ws.Range("M132").Value = -4053.2855
$ws.Range("N132").Value = -17605.25
$ws.Range("H134").Value = 29534.691
$ws.Range("J134").Value = 29534.691
$ws.Range("L134").Value = 88604.073
$ws.Range("N134").Value = -93674.073
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 632982.4
$ws.Range("I46").Value = 238.625
$ws.Range("J46").Value = 1265726.1
$ws.Range("K46").Value = 238.625
$ws.Range("L46").Value = 1265726.1
$ws.Range("M46").Value = -50.625
$ws.Range("N46").Value = -1266102.1
$ws.Range("H61").Value = 1804.6666
$ws.Range("I61").Value = 1931.5
$ws.Range("J61").Value = 1551
$ws.Range("K61").Value = 1931.5
$ws.Range("L61").Value = 1551
$ws.Range("M61").Value = -1729.5
$ws.Range("N61").Value = -1955
$ws.Range("H100").Value = 2462.6428
$ws.Range("I100").Value = 2227.5715
$ws.Range("J100").Value = 2697.7144
$ws.Range("K100").Value = 2227.5715
$ws.Range("L100").Value = 2697.7144
$ws.Range("M100").Value = -1686.5715
$ws.Range("N100").Value = -3779.7144
$ws.Range("H113").Value = 1804.6666
$ws.Range("I113").Value = 1931.5
$ws.Range("J113").Value = 1551
$ws.Range("K113").Value = 1931.5
$ws.Range("L113").Value = 1551
$ws.Range("M113").Value = 238.5
$ws.Range("N113").Value = -5891
$ws.Range("H132").Value = 2954.5
$ws.Range("I132").Value = 3046.5
$ws.Range("J132").Value = 2494.5
$ws.Range("K132").Value = 9139.5
$ws.Range("L132").Value = 7483.5
$ws.Range("M132").Value = -6609.5
$ws.Range("N132").Value = -12543.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 48253.906
$ws.Range("I100").Value = 77340.92
$ws.Range("J100").Value = 987.5
$ws.Range("K100").Value = 154681.84
$ws.Range("L100").Value = 1975
$ws.Range("M100").Value = -154140.84
$ws.Range("N100").Value = -3057
$ws.Range("H107").Value = 200575.6
$ws.Range("J107").Value = 333793
$ws.Range("L107").Value = 1001379
$ws.Range("N107").Value = -1005219
$ws.Range("H113").Value = 742.8182
$ws.Range("I113").Value = 512.1667
$ws.Range("J113").Value = 1019.6
$ws.Range("K113").Value = 1536.5001
$ws.Range("L113").Value = 3058.8
$ws.Range("M113").Value = 633.4999
$ws.Range("N113").Value = -7398.8
$ws.Range("H122").Value = 1632
$ws.Range("I122").Value = 1270.25
$ws.Range("K122").Value = 3810.75
$ws.Range("M122").Value = -1360.75
$ws.Range("H126").Value = 2201.2856
$ws.Range("I126").Value = 2639.5
$ws.Range("J126").Value = 2026
$ws.Range("K126").Value = 7918.5
$ws.Range("L126").Value = 6078
$ws.Range("M126").Value = -5448.5
$ws.Range("N126").Value = -11018
